$wb = $excel.ActiveWorkbook

# Sheet "NBR" - new Reaction_number values for rows 2..20 (column C)
$nbr = $wb.Worksheets.Item("NBR")
$nbrValues = @(498, 488, 476, 468, 456, 447, 447, 439, 433, 435, 431, 428, 422, 425, 439, 409, 410, 408, 406)
for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $nbr.Cells.Item($row, 3).Value = $nbrValues[$i]
}

# Sheet "BAR" - new Reaction_number values for rows 2..20 (column C)
$bar = $wb.Worksheets.Item("BAR")
$barValues = @(587, 593, 582, 585, 584, 585, 586, 583, 581, 581, 581, 581, 580, 576, 580, 579, 577, 580, 580)
for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $bar.Cells.Item($row, 3).Value = $barValues[$i]
}
